# Update the "想去人数" (number of people who want to go) counts in
# both the "展览" sheet and the "全部类型" sheet for the events whose
# F-column totals increased.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F8").Value = 2696
$wsExhibit.Range("F15").Value = 141
$wsExhibit.Range("F21").Value = 6067
$wsExhibit.Range("F23").Value = 1234
$wsExhibit.Range("F24").Value = 132
$wsExhibit.Range("F29").Value = 52
$wsExhibit.Range("F39").Value = 202

# --- Sheet "全部类型" (All types) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F11").Value = 2696
$wsAll.Range("F19").Value = 141
$wsAll.Range("F25").Value = 6067
$wsAll.Range("F27").Value = 1234
$wsAll.Range("F28").Value = 132
$wsAll.Range("F33").Value = 52
$wsAll.Range("F43").Value = 202
